# Apply the changes described by the diff to the workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: updated publish date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 used to be a duplicate "Contact" / "No display for ContactDetail" row;
# it becomes the new "Jurisdiction" / "United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was the other duplicate "Contact" row - remove it entirely, shifting
# everything below up by one row (dimension goes from A1:B21 to A1:B20).
$meta.Rows.Item(11).Delete()

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the "Extension" element) Short/Definition columns (K/L) change from
# the generic "Extension" / "An Extension" to the scoring-method specific text.
$elements.Range("K2").Value = "Scoring Method"
$elements.Range("L2").Value = "Method used to assess score the insight"
